$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 133.63637
$ws.Range("I8").Value = 133.63637
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 400.90911
$ws.Range("L8").Value = 0
$ws.Range("N8").Value = -261.90911
$ws.Range("M8").ClearContents()
# Row 13
$ws.Range("H13").Value = 10000
$ws.Range("J13").Value = 10000
$ws.Range("L13").Value = 10000
$ws.Range("N13").Value = -10338
# Row 15
$ws.Range("H15").Value = 1233.6833
$ws.Range("I15").Value = 1233.6833
$ws.Range("K15").Value = 3701.0499
$ws.Range("M15").Value = -3532.0499
# Row 33
$ws.Range("H33").Value = 896.1111
$ws.Range("I33").Value = 291.54166
$ws.Range("K33").Value = 291.54166
$ws.Range("M33").Value = -62.54165999999998
# Row 116
$ws.Range("H116").Value = 2159.0908
$ws.Range("I116").Value = 1630.5
$ws.Range("K116").Value = 1630.5
$ws.Range("M116").Value = 1811.5
# Row 129
$ws.Range("H129").Value = 2680.8276
$ws.Range("I129").Value = 11690.667
$ws.Range("J129").Value = 1025.9592
$ws.Range("K129").Value = 35072.001
$ws.Range("L129").Value = 3077.8776
$ws.Range("M129").Value = -30072.001
$ws.Range("N129").Value = -13077.8776
# Row 132
$ws.Range("H132").Value = 6950968.5
$ws.Range("I132").Value = 7149339
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 21448017
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -21445487
$ws.Range("N132").Value = -29060
# Row 135
$ws.Range("H135").Value = 965.3889
$ws.Range("I135").Value = 598.875
$ws.Range("J135").Value = 3897.5
$ws.Range("K135").Value = 5389.875
$ws.Range("L135").Value = 35077.5
$ws.Range("M135").Value = -2854.875
$ws.Range("N135").Value = -40147.5
# Row 137
$ws.Range("H137").Value = 2117.0386
$ws.Range("I137").Value = 1481.3684
$ws.Range("J137").Value = 3842.4285
$ws.Range("K137").Value = 4444.1052
$ws.Range("L137").Value = 11527.2855
$ws.Range("M137").Value = -1894.1052
$ws.Range("N137").Value = -16627.2855
# Row 138
$ws.Range("H138").Value = 6132.7964
$ws.Range("I138").Value = 1186.6842
$ws.Range("J138").Value = 17879.812
$ws.Range("K138").Value = 3560.0526
$ws.Range("L138").Value = 53639.436
$ws.Range("M138").Value = 1579.9474
$ws.Range("N138").Value = -63919.436

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2066.8235
$ws.Range("I61").Value = 1348.4445
$ws.Range("J61").Value = 2875
$ws.Range("K61").Value = 1348.4445
$ws.Range("L61").Value = 2875
$ws.Range("M61").Value = -1136.4445
$ws.Range("N61").Value = -3299
# Row 136
$ws.Range("H136").Value = 2066.8235
$ws.Range("I136").Value = 1348.4445
$ws.Range("J136").Value = 2875
$ws.Range("K136").Value = 4045.3335
$ws.Range("L136").Value = 8625
$ws.Range("M136").Value = -1495.3335
$ws.Range("N136").Value = -13725

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 400.66666
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 402
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 402
$ws.Range("M22").Value = -227
$ws.Range("N22").Value = -748
# Row 99
$ws.Range("H99").Value = 1544.2069
$ws.Range("I99").Value = 1372.8948
$ws.Range("J99").Value = 1869.7
$ws.Range("K99").Value = 1372.8948
$ws.Range("L99").Value = 1869.7
$ws.Range("M99").Value = 125.1052
$ws.Range("N99").Value = -4865.7
# Row 107
$ws.Range("H107").Value = 111208840
$ws.Range("I107").Value = 142979860
$ws.Range("J107").Value = 10310
$ws.Range("K107").Value = 142979860
$ws.Range("L107").Value = 10310
$ws.Range("M107").Value = -142977940
$ws.Range("N107").Value = -14150
# Row 134
$ws.Range("H134").Value = 2280.3333
$ws.Range("I134").Value = 2002.9788
$ws.Range("J134").Value = 4142.5713
$ws.Range("K134").Value = 6008.936400000001
$ws.Range("L134").Value = 12427.7139
$ws.Range("M134").Value = -3473.936400000001
$ws.Range("N134").Value = -17497.7139

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 19
$ws.Range("H19").Value = 20664.715
$ws.Range("I19").Value = 38.25
$ws.Range("K19").Value = 38.25
$ws.Range("M19").Value = 131.75
# Row 24
$ws.Range("H24").Value = 20664.715
$ws.Range("I24").Value = 38.25
$ws.Range("K24").Value = 38.25
$ws.Range("M24").Value = 131.75
# Row 31
$ws.Range("H31").Value = 27995.117
$ws.Range("I31").Value = 947.4
$ws.Range("J31").Value = 47314.914
$ws.Range("K31").Value = 947.4
$ws.Range("L31").Value = 47314.914
$ws.Range("M31").Value = -652.4
$ws.Range("N31").Value = -47904.914
# Row 34
$ws.Range("H34").Value = 27995.117
$ws.Range("I34").Value = 947.4
$ws.Range("J34").Value = 47314.914
$ws.Range("K34").Value = 947.4
$ws.Range("L34").Value = 47314.914
$ws.Range("M34").Value = -745.4
$ws.Range("N34").Value = -47718.914
# Row 58
$ws.Range("H58").Value = 3473.0352
$ws.Range("I58").Value = 797.2545
$ws.Range("J58").Value = 77057
$ws.Range("K58").Value = 797.2545
$ws.Range("L58").Value = 77057
$ws.Range("M58").Value = -594.2545
$ws.Range("N58").Value = -77463
# Row 68
$ws.Range("H68").Value = 18356.143
$ws.Range("J68").Value = 18356.143
$ws.Range("L68").Value = 18356.143
$ws.Range("N68").Value = -19854.143
# Row 71
$ws.Range("H71").Value = 18356.143
$ws.Range("J71").Value = 18356.143
$ws.Range("L71").Value = 55068.429
$ws.Range("N71").Value = -62556.429
# Row 132
$ws.Range("H132").Value = 37503804
$ws.Range("I132").Value = 34486310
$ws.Range("J132").Value = 45459000
$ws.Range("K132").Value = 103458930
$ws.Range("L132").Value = 136377000
$ws.Range("M132").Value = -103456400
$ws.Range("N132").Value = -136382060
# Row 134
$ws.Range("H134").Value = 1291.4
$ws.Range("I134").Value = 1160.9333
$ws.Range("J134").Value = 1682.8
$ws.Range("K134").Value = 3482.7999
$ws.Range("L134").Value = 5048.4
$ws.Range("M134").Value = -947.7999
$ws.Range("N134").Value = -10118.4
# Row 136
$ws.Range("H136").Value = 3473.0352
$ws.Range("I136").Value = 797.2545
$ws.Range("J136").Value = 77057
$ws.Range("K136").Value = 2391.7635
$ws.Range("L136").Value = 231171
$ws.Range("M136").Value = 158.2365
$ws.Range("N136").Value = -236271

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 10
$ws.Range("H10").Value = 189.4
$ws.Range("I10").Value = 68
$ws.Range("J10").Value = 675
$ws.Range("K10").Value = 204
$ws.Range("L10").Value = 2025
$ws.Range("M10").Value = -65
$ws.Range("N10").Value = -2303

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 3176.2
$ws.Range("I132").Value = 2087.8948
$ws.Range("J132").Value = 4468.5625
$ws.Range("K132").Value = 6263.6844
$ws.Range("L132").Value = 13405.6875
$ws.Range("M132").Value = -3733.6844
$ws.Range("N132").Value = -18465.6875

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1355.7084
$ws.Range("I22").Value = 1667
$ws.Range("J22").Value = 1092.3077
$ws.Range("K22").Value = 1667
$ws.Range("L22").Value = 1092.3077
$ws.Range("M22").Value = -1372
$ws.Range("N22").Value = -1682.3077
# Row 27
$ws.Range("H27").Value = 1355.7084
$ws.Range("I27").Value = 1667
$ws.Range("J27").Value = 1092.3077
$ws.Range("K27").Value = 1667
$ws.Range("L27").Value = 1667
$ws.Range("M27").Value = -1560
$ws.Range("N27").Value = -1306.3077
# Row 132
$ws.Range("H132").Value = 3646.6562
$ws.Range("I132").Value = 3507.2593
$ws.Range("J132").Value = 4399.4
$ws.Range("K132").Value = 10521.7779
$ws.Range("L132").Value = 13198.2
$ws.Range("M132").Value = -7991.777900000001
$ws.Range("N132").Value = -18258.2
# Row 136
$ws.Range("H136").Value = 1794.45
$ws.Range("I136").Value = 1397.4286
$ws.Range("K136").Value = 4192.2858
$ws.Range("M136").Value = -1642.2858

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2334.689
$ws.Range("I132").Value = 2144.8333
$ws.Range("K132").Value = 6434.499899999999
$ws.Range("M132").Value = -3904.499899999999
# Row 136
$ws.Range("H136").Value = 751.9268
$ws.Range("I136").Value = 439.84848
$ws.Range("K136").Value = 1319.54544
$ws.Range("M136").Value = 1230.45456

Write-Output "Applied Aegis_Profits updates"